$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.026688194879163
$ws.Cells.Item(2, 4).Value = 1.035842992883215
$ws.Cells.Item(2, 5).Value = 1.026877410206643
$ws.Cells.Item(2, 6).Value = 1.046628866459433
$ws.Cells.Item(2, 9).Value = 1.033561569800667
$ws.Cells.Item(2, 10).Value = 1.0318502428324
$ws.Cells.Item(2, 11).Value = 1.038638437759446
$ws.Cells.Item(2, 12).Value = 1.029698774358945
$ws.Cells.Item(2, 13).Value = 1.049393768197427
$ws.Cells.Item(2, 14).Value = 1.033315588532099

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.027616856318629
$ws.Cells.Item(3, 4).Value = 1.036535536090016
$ws.Cells.Item(3, 5).Value = 1.027665520847968
$ws.Cells.Item(3, 6).Value = 1.04747817525519
$ws.Cells.Item(3, 9).Value = 1.033705487522319
$ws.Cells.Item(3, 10).Value = 1.032418877727394
$ws.Cells.Item(3, 11).Value = 1.039140770230073
$ws.Cells.Item(3, 12).Value = 1.030294544744671
$ws.Cells.Item(3, 13).Value = 1.05005466028125
$ws.Cells.Item(3, 14).Value = 1.033885030953868

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.028218386999037
$ws.Cells.Item(4, 4).Value = 1.036984016336785
$ws.Cells.Item(4, 5).Value = 1.028176398268278
$ws.Cells.Item(4, 6).Value = 1.048028420948584
$ws.Cells.Item(4, 9).Value = 1.033797429730266
$ws.Cells.Item(4, 10).Value = 1.032786839761339
$ws.Cells.Item(4, 11).Value = 1.039465495102618
$ws.Cells.Item(4, 12).Value = 1.030680310402235
$ws.Cells.Item(4, 13).Value = 1.05048234280247
$ws.Cells.Item(4, 14).Value = 1.034253515536107

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.028471418845127
$ws.Cells.Item(5, 4).Value = 1.037172641652454
$ws.Cells.Item(5, 5).Value = 1.028391389162149
$ws.Cells.Item(5, 6).Value = 1.04825990675641
$ws.Cells.Item(5, 9).Value = 1.033835798607288
$ws.Cells.Item(5, 10).Value = 1.032941534004355
$ws.Cells.Item(5, 11).Value = 1.039601932155528
$ws.Cells.Item(5, 12).Value = 1.030842547904812
$ws.Cells.Item(5, 13).Value = 1.050662148913899
$ws.Cells.Item(5, 14).Value = 1.034408429462692

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.028513912670193
$ws.Cells.Item(6, 4).Value = 1.037204317565054
$ws.Cells.Item(6, 5).Value = 1.028427499819698
$ws.Cells.Item(6, 6).Value = 1.048298783721755
$ws.Cells.Item(6, 9).Value = 1.033842224260809
$ws.Cells.Item(6, 10).Value = 1.032967507992006
$ws.Cells.Item(6, 11).Value = 1.039624835970373
$ws.Cells.Item(6, 12).Value = 1.030869791874101
$ws.Cells.Item(6, 13).Value = 1.050692339599779
$ws.Cells.Item(6, 14).Value = 1.034434440336384

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.028221767441524
$ws.Cells.Item(7, 4).Value = 1.036986536428572
$ws.Cells.Item(7, 5).Value = 1.028179270132361
$ws.Cells.Item(7, 6).Value = 1.048031513436875
$ws.Cells.Item(7, 9).Value = 1.033797943532301
$ws.Cells.Item(7, 10).Value = 1.032788906782989
$ws.Cells.Item(7, 11).Value = 1.039467318486368
$ws.Cells.Item(7, 12).Value = 1.030682477987058
$ws.Cells.Item(7, 13).Value = 1.050484745349454
$ws.Cells.Item(7, 14).Value = 1.034255585493164

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.027001910309795
$ws.Cells.Item(8, 4).Value = 1.036076965913624
$ws.Cells.Item(8, 5).Value = 1.027143565064206
$ws.Cells.Item(8, 6).Value = 1.046915751234681
$ws.Cells.Item(8, 9).Value = 1.033610451879619
$ws.Cells.Item(8, 10).Value = 1.032042411414565
$ws.Cells.Item(8, 11).Value = 1.038808268459098
$ws.Cells.Item(8, 12).Value = 1.029900062704631
$ws.Cells.Item(8, 13).Value = 1.049617110365978
$ws.Cells.Item(8, 14).Value = 1.033508030015693

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.024857202767679
$ws.Cells.Item(9, 4).Value = 1.034477002737008
$ws.Cells.Item(9, 5).Value = 1.02532561252346
$ws.Cells.Item(9, 6).Value = 1.044954967267618
$ws.Cells.Item(9, 9).Value = 1.033271040539908
$ws.Cells.Item(9, 10).Value = 1.030727172751829
$ws.Cells.Item(9, 11).Value = 1.037644556744857
$ws.Cells.Item(9, 12).Value = 1.028523409290068
$ws.Cells.Item(9, 13).Value = 1.048088606869762
$ws.Cells.Item(9, 14).Value = 1.032190923563195

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.023430719688471
$ws.Cells.Item(10, 4).Value = 1.033412349253733
$ws.Cells.Item(10, 5).Value = 1.024118503823213
$ws.Cells.Item(10, 6).Value = 1.043651465828218
$ws.Cells.Item(10, 9).Value = 1.033038733283272
$ws.Cells.Item(10, 10).Value = 1.029850535482342
$ws.Cells.Item(10, 11).Value = 1.036867218747387
$ws.Cells.Item(10, 12).Value = 1.027607097716884
$ws.Cells.Item(10, 13).Value = 1.047069947533065
$ws.Cells.Item(10, 14).Value = 1.031313041368234

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.022813839180641
$ws.Cells.Item(11, 4).Value = 1.032951834208524
$ws.Cells.Item(11, 5).Value = 1.023596984248468
$ws.Cells.Item(11, 6).Value = 1.043087931692967
$ws.Cells.Item(11, 9).Value = 1.032936718312604
$ws.Cells.Item(11, 10).Value = 1.029471000460113
$ws.Cells.Item(11, 11).Value = 1.036530274726925
$ws.Cells.Item(11, 12).Value = 1.027210685186423
$ws.Cells.Item(11, 13).Value = 1.046628956470216
$ws.Cells.Item(11, 14).Value = 1.030932967362741

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.022584823026023
$ws.Cells.Item(12, 4).Value = 1.032780853177119
$ws.Cells.Item(12, 5).Value = 1.023403445445777
$ws.Cells.Item(12, 6).Value = 1.042878745561604
$ws.Cells.Item(12, 9).Value = 1.032898612000219
$ws.Cells.Item(12, 10).Value = 1.02933003357192
$ws.Cells.Item(12, 11).Value = 1.036405067079569
$ws.Cells.Item(12, 12).Value = 1.02706349465339
$ws.Cells.Item(12, 13).Value = 1.046465168510864
$ws.Cells.Item(12, 14).Value = 1.030791800285397

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.022633942286989
$ws.Cells.Item(13, 4).Value = 1.032817525779945
$ws.Cells.Item(13, 5).Value = 1.02344495215117
$ws.Cells.Item(13, 6).Value = 1.042923610535625
$ws.Cells.Item(13, 9).Value = 1.032906795581361
$ws.Cells.Item(13, 10).Value = 1.029360271008641
$ws.Cells.Item(13, 11).Value = 1.03643192686265
$ws.Cells.Item(13, 12).Value = 1.027095065026701
$ws.Cells.Item(13, 13).Value = 1.046500300847662
$ws.Cells.Item(13, 14).Value = 1.030822080662746

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.02279490616439
$ws.Cells.Item(14, 4).Value = 1.032937699328762
$ws.Cells.Item(14, 5).Value = 1.023580982655331
$ws.Cells.Item(14, 6).Value = 1.043070637527117
$ws.Cells.Item(14, 9).Value = 1.032933572783023
$ws.Cells.Item(14, 10).Value = 1.029459347896715
$ws.Cells.Item(14, 11).Value = 1.03651992607027
$ws.Cells.Item(14, 12).Value = 1.027198517247713
$ws.Cells.Item(14, 13).Value = 1.046615417376523
$ws.Cells.Item(14, 14).Value = 1.030921298251366

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.022894097353895
$ws.Cells.Item(15, 4).Value = 1.033011752172896
$ws.Cells.Item(15, 5).Value = 1.02366481901073
$ws.Cells.Item(15, 6).Value = 1.043161243721106
$ws.Cells.Item(15, 9).Value = 1.032950042839471
$ws.Cells.Item(15, 10).Value = 1.029520393705027
$ws.Cells.Item(15, 11).Value = 1.036574138479457
$ws.Cells.Item(15, 12).Value = 1.027262264863186
$ws.Cells.Item(15, 13).Value = 1.046686346600874
$ws.Cells.Item(15, 14).Value = 1.030982430751729

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.02347167688499
$ws.Cells.Item(16, 4).Value = 1.033442922516784
$ws.Cells.Item(16, 5).Value = 1.024153140067272
$ws.Cells.Item(16, 6).Value = 1.043688884644774
$ws.Cells.Item(16, 9).Value = 1.033045473720788
$ws.Cells.Item(16, 10).Value = 1.029875725218122
$ws.Cells.Item(16, 11).Value = 1.036889573306449
$ws.Cells.Item(16, 12).Value = 1.027633413921366
$ws.Cells.Item(16, 13).Value = 1.047099216786988
$ws.Cells.Item(16, 14).Value = 1.031338266876328

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.023834191194914
$ws.Cells.Item(17, 4).Value = 1.033713515782109
$ws.Cells.Item(17, 5).Value = 1.024459764346626
$ws.Cells.Item(17, 6).Value = 1.044020099518387
$ws.Cells.Item(17, 9).Value = 1.033104954055828
$ws.Cells.Item(17, 10).Value = 1.030098630709458
$ws.Cells.Item(17, 11).Value = 1.037087343955233
$ws.Cells.Item(17, 12).Value = 1.027866322179902
$ws.Cells.Item(17, 13).Value = 1.047358225853118
$ws.Cells.Item(17, 14).Value = 1.03156148891903

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.02404571639258
$ws.Cells.Item(18, 4).Value = 1.033871395130744
$ws.Cells.Item(18, 5).Value = 1.024638725544083
$ws.Cells.Item(18, 6).Value = 1.044213377362984
$ws.Cells.Item(18, 9).Value = 1.033139510455822
$ws.Cells.Item(18, 10).Value = 1.030228652915278
$ws.Cells.Item(18, 11).Value = 1.037202666231269
$ws.Cells.Item(18, 12).Value = 1.028002207937411
$ws.Cells.Item(18, 13).Value = 1.04750931054211
$ws.Cells.Item(18, 14).Value = 1.031691695771295

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.024117853999769
$ws.Cells.Item(19, 4).Value = 1.033925235810669
$ws.Cells.Item(19, 5).Value = 1.024699765713041
$ws.Cells.Item(19, 6).Value = 1.0442792946157
$ws.Cells.Item(19, 9).Value = 1.033151269959471
$ws.Cells.Item(19, 10).Value = 1.030272987976631
$ws.Cells.Item(19, 11).Value = 1.037241982344796
$ws.Cells.Item(19, 12).Value = 1.028048547252877
$ws.Cells.Item(19, 13).Value = 1.047560828012834
$ws.Cells.Item(19, 14).Value = 1.031736093793521

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.023795288895121
$ws.Cells.Item(20, 4).Value = 1.033684478825516
$ws.Cells.Item(20, 5).Value = 1.024426854821678
$ws.Cells.Item(20, 6).Value = 1.043984554408495
$ws.Cells.Item(20, 9).Value = 1.033098586594028
$ws.Cells.Item(20, 10).Value = 1.030074714513571
$ws.Cells.Item(20, 11).Value = 1.037066128544291
$ws.Cells.Item(20, 12).Value = 1.027841329770514
$ws.Cells.Item(20, 13).Value = 1.047330435681526
$ws.Cells.Item(20, 14).Value = 1.031537538759401

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.022747502972496
$ws.Cells.Item(21, 4).Value = 1.032902309137246
$ws.Cells.Item(21, 5).Value = 1.023540920169739
$ws.Cells.Item(21, 6).Value = 1.043027337987219
$ws.Cells.Item(21, 9).Value = 1.032925693450961
$ws.Cells.Item(21, 10).Value = 1.029430171956548
$ws.Cells.Item(21, 11).Value = 1.036494013912879
$ws.Cells.Item(21, 12).Value = 1.027168051626766
$ws.Cells.Item(21, 13).Value = 1.046581517981287
$ws.Cells.Item(21, 14).Value = 1.030892080878018

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.022089417461721
$ws.Cells.Item(22, 4).Value = 1.032410961345865
$ws.Cells.Item(22, 5).Value = 1.022984921494965
$ws.Cells.Item(22, 6).Value = 1.042426283071127
$ws.Cells.Item(22, 9).Value = 1.032815754237666
$ws.Cells.Item(22, 10).Value = 1.029024976924549
$ws.Cells.Item(22, 11).Value = 1.036134004572299
$ws.Cells.Item(22, 12).Value = 1.026745052177762
$ws.Cells.Item(22, 13).Value = 1.046110735584374
$ws.Cells.Item(22, 14).Value = 1.030486310422596

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.022438214398174
$ws.Cells.Item(23, 4).Value = 1.03267139249632
$ws.Cells.Item(23, 5).Value = 1.023279569372081
$ws.Cells.Item(23, 6).Value = 1.042744838661869
$ws.Cells.Item(23, 9).Value = 1.032874151907544
$ws.Cells.Item(23, 10).Value = 1.029239772948323
$ws.Cells.Item(23, 11).Value = 1.036324880204213
$ws.Cells.Item(23, 12).Value = 1.02696926172345
$ws.Cells.Item(23, 13).Value = 1.046360297042895
$ws.Cells.Item(23, 14).Value = 1.030701411481362

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.023812866932866
$ws.Cells.Item(24, 4).Value = 1.033697599231669
$ws.Cells.Item(24, 5).Value = 1.024441724872175
$ws.Cells.Item(24, 6).Value = 1.044000615448684
$ws.Cells.Item(24, 9).Value = 1.033101464200762
$ws.Cells.Item(24, 10).Value = 1.030085521197066
$ws.Cells.Item(24, 11).Value = 1.037075714980462
$ws.Cells.Item(24, 12).Value = 1.027852622658644
$ws.Cells.Item(24, 13).Value = 1.047342992835951
$ws.Cells.Item(24, 14).Value = 1.031548360789626

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.025411080710571
$ws.Cells.Item(25, 4).Value = 1.03489028801796
$ws.Cells.Item(25, 5).Value = 1.025794746968067
$ws.Cells.Item(25, 6).Value = 1.045461234466607
$ws.Cells.Item(25, 9).Value = 1.033271040539908
$ws.Cells.Item(25, 10).Value = 1.031067165114503
$ws.Cells.Item(25, 11).Value = 1.037945679014957
$ws.Cells.Item(25, 12).Value = 1.028879055510103
$ws.Cells.Item(25, 13).Value = 1.048483707603579
$ws.Cells.Item(25, 14).Value = 1.032531398754022
